$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45904
$ws.Range("B2").Value = 91.23999999999999
$ws.Range("C2").Value = 85
$ws.Range("D2").Value = 79.91
$ws.Range("E2").Value = 72.27
$ws.Range("F2").Value = 70
$ws.Range("G2").Value = 70.7
$ws.Range("H2").Value = 81.02
$ws.Range("I2").Value = 103.55
$ws.Range("J2").Value = 93.53
$ws.Range("K2").Value = 50
$ws.Range("L2").Value = 18.07
$ws.Range("M2").Value = 5.79
$ws.Range("N2").Value = 3.2
$ws.Range("O2").Value = 2.01
$ws.Range("P2").Value = 1.6
$ws.Range("Q2").Value = 2.01
$ws.Range("R2").Value = 3.7
$ws.Range("S2").Value = 7.01
$ws.Range("T2").Value = 28.19
$ws.Range("U2").Value = 66.04000000000001
$ws.Range("V2").Value = 102.5
$ws.Range("W2").Value = 112.99
$ws.Range("X2").Value = 100.45
$ws.Range("Y2").Value = 85.5
$ws.Range("Z2").Value = 55.68
$ws.Range("AB2").Value = 100.36
$ws.Range("AD2").Value = 107.74
$ws.Range("AF2").Value = 92.98
